$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$v = $ws.Range("A1").Value2
Write-Output ($v.Length)
Write-Output ([string]::Join(",", ($v.ToCharArray() | ForEach-Object { [int]$_ })))
